$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# Insert a new worksheet row at row 35 by copying the row above it (row 34, "Chemi" /
# Healthcare), which both shifts everything below down by one row AND carries over the
# correct cell styling/formula for the new table row, the way Excel's own
# "insert copied cells" does. "Clini" sorts between "Chemi" and "Disab" (row 35 was "Disab").
$ws.Rows.Item(34).Copy()
$ws.Rows.Item(35).Insert()

# Grow the table definition (ref / autoFilter / sortState) to cover the new row.
$tbl.Resize($ws.Range("A1:C58"))

# Replace the duplicated content in the new row with the actual new keyword.
$ws.Range("A35").Value = "Clini"
$ws.Range("C35").Value = "Healthcare"

# Re-apply the calculated column formula across the whole table (the row that used to be
# the table's last row can otherwise end up with a stale/broken structured reference
# after the resize above).
$tbl.ListColumns.Item(2).DataBodyRange.Formula = "=LEN(Cluster_Keywords[[#This Row],[Stem]])"

# The "duplicate stem" highlighting rule covered A2:A16, A46:A57, A18:A22 and A25:A44; those
# row numbers need to follow the rows that shifted down by the insert above (-> A47:A58 and
# A25:A45). Recreate it as one rule per block, restyled to match the original red-on-pink
# duplicate-value look.
$ws.UsedRange.FormatConditions.Item(1).Delete()
foreach ($addr in @("A2:A16", "A18:A22", "A25:A45", "A47:A58")) {
    $rule = $ws.Range($addr).FormatConditions.AddUniqueValues()
    $rule.DupeUnique = 1
    $rule.Font.Color = 393372
    $rule.Interior.Color = 13551615
}

# Move the active selection to the new location, matching the saved view state.
[void]$ws.Range("A36").Select()
